# Update odds values in "Jogos da Semana" worksheet per upstream source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("K5").Value = 1.83
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.25
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 451

# Row 8
$ws.Range("G8").Value = 1.53
$ws.Range("H8").Value = 4.1
$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 2.1
$ws.Range("Z8").Value = 11
$ws.Range("AD8").Value = 7.5
$ws.Range("AE8").Value = 17
$ws.Range("AN8").Value = 3.5

# Row 10
$ws.Range("G10").Value = 3.4
$ws.Range("I10").Value = 2.2
$ws.Range("W10").Value = 9
$ws.Range("AK10").Value = 19
$ws.Range("AW10").Value = 4
